# "added Voci 4 #2"
# Adds the "Lehrbuch 4" vocabulary section: fills in the previously-blank
# column D (chapter label) for rows 396-470, fixes a typo in B461
# ("die Verwanten" -> "die Verwandten"), and appends 10 new vocabulary
# rows (471-480) for the new "Lehrbuch 4" chapter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 396-480: (row, ColumnA, ColumnB, ColumnC). Column D is always
# "Lehrbuch 4" for all of these rows.
$rows = @(
    @(396, 'το δέντρο', 'der Baum', 'Nomen'),
    @(397, 'το γενεαλογικό δέντρο', 'der Stammbaum', 'Nomen'),
    @(398, 'ο παππούς', 'der Grossvater', 'Nomen'),
    @(399, 'η γιαγιά', 'die Grossmutter', 'Nomen'),
    @(400, 'η μητέρα', 'die Mutter', 'Nomen'),
    @(401, 'ο πατέρας', 'der Vater', 'Nomen'),
    @(402, 'η νύφη', 'die Schwägerin', 'Nomen'),
    @(403, 'ο αδελφός', 'der Bruder', 'Nomen'),
    @(404, 'η αδελφή', 'die Schwester', 'Nomen'),
    @(405, 'ο γαμπρός', 'der Schwager', 'Nomen'),
    @(406, 'η γαμπρή', 'die Schwägerin', 'Nomen'),
    @(407, 'ο γιος', 'der Sohn', 'Nomen'),
    @(408, 'πόσος, -η, -ο', 'wie viele', 'Adjektiv'),
    @(409, 'το κλαδί', 'der Ast, Zweig', 'Nomen'),
    @(410, 'δικός, -ιά, -ή, -ό', 'mein', 'Pronomen'),
    @(411, 'το κορίτσι', 'das Mädchen', 'Nomen'),
    @(412, 'εκεί', 'da, dort', 'Partikel'),
    @(413, 'η οικογένεια', 'die Familie', 'Nomen'),
    @(414, 'το άλμπουμ', 'das Album', 'Nomen'),
    @(415, 'ποιος, -α, -ο', 'wer', 'Pronomen'),
    @(416, 'πρώτος, -η, -ο', 'erster', 'Adjektiv'),
    @(417, 'η σελίδα', 'die Seite', 'Nomen'),
    @(418, 'ίδιος, -η, -ο', 'gleich (wie)', 'Adjektiv'),
    @(419, 'το όνομα', 'der Name', 'Nomen'),
    @(420, 'παλιός, -ή, -ό', 'alt', 'Adjektiv'),
    @(421, 'ασπρόμαυρος, -η, -ο', 'schwarzweiss', 'Adjektiv'),
    @(422, 'έγχρωμος, -η, -ο', 'farbig', 'Adjektiv'),
    @(423, 'οικογενειακός, -ή, -ό', 'Familien-', 'Adjektiv'),
    @(424, 'ο γάμος', 'die Hochzeit', 'Nomen'),
    @(425, 'τα εγγόνια', 'die Enkel', 'Nomen'),
    @(426, 'πόπο', 'Oh', 'Partikel'),
    @(427, 'η θεία', 'die Tante', 'Nomen'),
    @(428, 'ο θείος', 'der Onkel', 'Nomen'),
    @(429, 'από τρία', 'ab drei (Anzahl)', 'Spruch'),
    @(430, 'παντρεμένος, -η, -ο', 'verheiratet', 'Adjektiv'),
    @(431, 'άρα', 'also', 'Partikel'),
    @(432, 'τα ξαδέλφια', 'die Cousins & Cousinen', 'Nomen'),
    @(433, 'σημαίνω', 'bedeuten', 'Verb'),
    @(434, 'ότι', 'dass', 'Partikel'),
    @(435, 'η ξαδέλφη', 'die Cousine', 'Nomen'),
    @(436, 'ο ξάδελφος', 'der Cousin', 'Nomen'),
    @(437, 'υπάρχει', 'es gibt', 'Verb'),
    @(438, 'τα αδέλφια', 'die Geschwister', 'Nomen'),
    @(439, 'το γούστο', 'der Geschmack', 'Nomen'),
    @(440, 'η σχέση', 'die Beziehung', 'Nomen'),
    @(441, 'έχω σχέση με', 'etw. zu tun haben mit', 'Spruch'),
    @(442, 'συνήθως', 'gewöhnlich', 'Partikel'),
    @(443, 'το μπέρδεμα', 'das Durcheinander', 'Nomen'),
    @(444, 'η ανιψιά', 'die Nichte', 'Nomen'),
    @(445, 'η μαμά', 'die Mama', 'Nomen'),
    @(446, 'ο μπαμπάς', 'der Papa', 'Nomen'),
    @(447, 'επίσης', 'ebenfalls', 'Partikel'),
    @(448, 'εύκολος, -η, -ο | απλός, -ή, -ό', 'einfach', 'Adjektiv'),
    @(449, 'διπλός, -ή, -ό', 'doppelt', 'Adjektiv'),
    @(450, 'τριπλός, -ή, -ό', 'dreifach', 'Adjektiv'),
    @(451, 'μοναδικός, -ή, -ό', 'einzigartig', 'Adjektiv'),
    @(452, 'μόνος, -η, -ο', 'einzig', 'Adjektiv'),
    @(453, 'κι', 'und, auch', 'Partikel'),
    @(454, 'δίνω', 'geben', 'Verb'),
    @(455, 'η αλήθεια', 'die Wahrheit', 'Nomen'),
    @(456, 'αυτό ειναι αλήθεια', 'das stimmt', 'Spruch'),
    @(457, 'το μοναχοπαίδι', 'das Einzelkind', 'Nomen'),
    @(458, 'η μοναχοκόρη', 'das Einzelkind (f)', 'Nomen'),
    @(459, 'ο μοναχογιός', 'das Einzelkind (m)', 'Nomen'),
    @(460, 'πότε', 'wann', 'Partikel'),
    @(461, 'οι συγγενείς', 'die Verwandten', 'Nomen'),
    @(462, 'οι γονείς', 'die Eltern', 'Nomen'),
    @(463, 'τα παιδιά', 'die Kinder', 'Nomen'),
    @(464, 'τα ανίψια', 'die Neffen und Nichten', 'Nomen'),
    @(465, 'τα πεθερικά', 'die Schwiegereltern', 'Nomen'),
    @(466, 'το ζευγάρι', 'das Paar', 'Nomen'),
    @(467, 'οι σύζυγοι', 'das Ehepaar', 'Nomen'),
    @(468, 'ο αρραβώνας', 'die Verlobung', 'Nomen'),
    @(469, 'ο χωρισμός', 'die Trennung', 'Nomen'),
    @(470, 'το διαζύγιο', 'die Schweidung', 'Nomen'),
    @(471, 'ο άντρας', 'der Mann', 'Nomen'),
    @(472, 'ο σύζυγος', 'der Ehemann', 'Nomen'),
    @(473, 'ανύπαντρος, ελεύθερος', 'ledig', 'Adjektiv'),
    @(474, 'ο αρραβωνιαστικός', 'der Verlobte', 'Nomen'),
    @(475, 'η αρραβωνιαστικά', 'die Verlobte', 'Nomen'),
    @(476, 'χωρισμένος', 'getrennt', 'Adjektiv'),
    @(477, 'διαζευγμένος', 'geschieden', 'Adjektiv'),
    @(478, 'ο χήρος', 'der Witwer', 'Nomen'),
    @(479, 'η χήρα', 'die Witwe', 'Nomen'),
    @(480, 'η σύζυγη', 'die Ehefrau', 'Nomen')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = 'Lehrbuch 4'
}

# Column A is widened slightly (20.796875 chars) to fit the new entries;
# columns B-D stay at the original 18.1328125. The COM layer here only
# supports ~1/6-character precision on ColumnWidth, so 19.92 is the
# closest input that lands on the nearest representable width.
$ws.Columns.Item(1).ColumnWidth = 19.92

# Scroll/selection state left by the author after typing the new rows.
$excel.ActiveWindow.ScrollRow = 457
[void]$ws.Range("D469:D480").Select()
